# Update the "Prix Spot" sheet with a new day column (22-jul) in column AM.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Header cell AM1: new date label, copy formatting from AL1 (the previous last header cell)
$ws.Range("AL1").Copy() | Out-Null
$ws.Range("AM1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("AM1").Value = "22-jul"

# New data values for column AM, rows 2-25
$values = @{
    2  = 81.13
    3  = 71.19
    4  = 67.28
    5  = 54.57
    6  = 46.63
    7  = 63.19
    8  = 79.06999999999999
    9  = 61.41
    10 = 76.26000000000001
    11 = 50.52
    12 = 28.86
    13 = 29.99
    14 = 27.67
    15 = 12.63
    16 = 4.18
    17 = 1.14
    18 = 27.05
    19 = 35.01
    20 = 69.06999999999999
    21 = 91.59
    22 = 113.92
    23 = 113.93
    24 = 109.1
    25 = 103.6
}

foreach ($row in $values.Keys) {
    $ws.Range("AM$row").Value = $values[$row]
}
